$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new value would otherwise be
# auto-parsed as a number by Excel (so they stay text, matching the source data).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated price (D) and volume (E) figures.
$ws.Range("D2").Value = "64.007.31"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "3.061.58"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "560.49"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "143.54"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.059.83"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  -10.67%  "
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  +8.90%  "
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "35.74"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "3.563.46"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "64.048.71"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "3.063.41"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "0.110"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "478.07"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "13.99"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").Value = "14.33"
$ws.Range("E24").Value = "  +9.55%  "
$ws.Range("D25").Value = "82.67"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "8.10"
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("D29").Value = "2.05"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "2.46"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").Value = "5.79"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "6.24"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "54.67"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "0.0412"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").Value = "451.87"
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").Value = "  +5.67%  "
$ws.Range("D41").Value = "3.030.35"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("D45").Value = "27.86"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  +8.11%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("D49").Value = "119.36"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").Value = "  +1.65%  "
